$d = $word.ActiveDocument

function Get-MatchRange($searchText) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $rng.Find.Found) {
        throw "Could not find text: $searchText"
    }
    return $d.Range($rng.Start, $rng.End)
}

function Apply-NoteEdit($oldFirstRunText, $newFirstRunText, $lastRunText) {
    # Locate the three runs that make up this "Note: ..." sentence:
    #   run1 = the sentence up to (but not including) the trailing " " + last word/phrase
    #   run2 = a single space " "
    #   run3 = the trailing word/phrase (e.g. "installation submissions" / "venues")
    $r1 = Get-MatchRange $oldFirstRunText
    $r2 = $d.Range($r1.End, $r1.End + 1)
    $r3 = $d.Range($r2.End, $r2.End + $lastRunText.Length)

    # Sanity-check the structure before mutating anything.
    if ($r2.Text -ne " ") { throw "Expected a single space run, got [$($r2.Text)]" }
    if ($r3.Text -ne $lastRunText) { throw "Expected [$lastRunText], got [$($r3.Text)]" }

    # Give the three runs temporarily distinct Bold/Italic combinations so that
    # the engine's "coalesce adjacent same-format runs" behaviour (triggered by
    # any text edit in the paragraph) can't merge them into one run while we
    # rewrite run1's text below. These combinations are deliberately chosen from
    # the same {Bold,Italic} values we want at the end, so no cleanup is needed -
    # the final formatting pass below converges everything to the real target.
    $r1.Font.Bold = 1
    $r1.Font.Italic = 0

    $r2.Font.Bold = 1
    $r2.Font.Italic = 1

    $r3.Font.Bold = 0
    $r3.Font.Italic = 1

    # Re-fetch run1 (formatting writes don't move boundaries, but stay safe) and
    # rewrite its text - this is now safe from merging neighbouring runs since
    # r1/r2/r3 no longer share identical formatting.
    $r1b = $d.Range($r1.Start, $r1.End)
    $r1b.Text = $newFirstRunText

    # Converge the whole sentence (all three runs) to the real target formatting:
    # italic, not bold. This is a pure formatting write (no text change) so it
    # does not trigger any run-merging - but it does leave all three runs with
    # identical <w:i/> formatting, exactly matching the desired end state.
    $whole = $d.Range($r1b.Start, $r3.End)
    $whole.Font.Bold = 0
    $whole.Font.Italic = 1
}

Apply-NoteEdit `
    "Note: a link to supporting material is required for performance/music or" `
    "note: a link to supporting material is required for performance/music or" `
    "installation submissions"

Apply-NoteEdit `
    "Note: all tech requirements are subject to available equipment and" `
    "note: all tech requirements are subject to available equipment and" `
    "venues"
